# Insert a new data row before the existing row 180 ("Acelga" price entry
# dated 2021-03-02 / serial 44257), shifting all rows from 180 through 268
# down by one (to 181 through 269). The new row carries a fresh weekly
# price observation (serial 44726 = 2022-06-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(180).Insert()

$ws.Range("A180").Value = 11
$ws.Range("B180").Value = "Vega Monumental Concepción"
$ws.Range("C180").Value = "Bíobío"
$ws.Range("D180").Value = 44726
$ws.Range("E180").Value = 8
$ws.Range("F180").Value = 100112009
$ws.Range("G180").Value = "Acelga"
$ws.Range("H180").Value = "Sin especificar"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 300
$ws.Range("K180").Value = 600
$ws.Range("L180").Value = 650
$ws.Range("M180").Value = 625
$ws.Range("N180").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O180").Value = "Región de Ñuble"
$ws.Range("P180").Value = 625
$ws.Range("Q180").Value = 1
$ws.Range("R180").Value = "Hortaliza"
